$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 857.55554
$ws.Cells.Item(2, 9).Value = 802.5714
$ws.Cells.Item(2, 10).Value = 1050
$ws.Cells.Item(2, 11).Value = 802.5714
$ws.Cells.Item(2, 12).Value = 1050
$ws.Cells.Item(2, 13).Value = -689.5714
$ws.Cells.Item(2, 14).Value = -1276

$ws.Cells.Item(18, 8).Value = 1425
$ws.Cells.Item(18, 9).Value = 1233.3334
$ws.Cells.Item(18, 10).Value = 2000
$ws.Cells.Item(18, 11).Value = 1233.3334
$ws.Cells.Item(18, 12).Value = 2000
$ws.Cells.Item(18, 13).Value = -949.3334
$ws.Cells.Item(18, 14).Value = -2568

$ws.Cells.Item(39, 8).Value = 304.41177
$ws.Cells.Item(39, 9).Value = 37.81818
$ws.Cells.Item(39, 10).Value = 793.1667
$ws.Cells.Item(39, 11).Value = 113.45454
$ws.Cells.Item(39, 12).Value = 2379.5001
$ws.Cells.Item(39, 13).Value = 182.54546
$ws.Cells.Item(39, 14).Value = -2971.5001

$ws.Cells.Item(43, 8).Value = 2499.5
$ws.Cells.Item(43, 9).Value = 2999
$ws.Cells.Item(43, 10).Value = 2000
$ws.Cells.Item(43, 11).Value = 2999
$ws.Cells.Item(43, 12).Value = 2000
$ws.Cells.Item(43, 13).Value = -2930

$ws.Cells.Item(64, 8).Value = 5248.36
$ws.Cells.Item(64, 9).Value = 5072.778
$ws.Cells.Item(64, 10).Value = 6828.6
$ws.Cells.Item(64, 11).Value = 5072.778
$ws.Cells.Item(64, 12).Value = 6828.6
$ws.Cells.Item(64, 13).Value = -4824.778
$ws.Cells.Item(64, 14).Value = -7324.6

$ws.Cells.Item(67, 8).Value = 5248.36
$ws.Cells.Item(67, 9).Value = 5072.778
$ws.Cells.Item(67, 10).Value = 6828.6
$ws.Cells.Item(67, 11).Value = 5072.778
$ws.Cells.Item(67, 12).Value = 6828.6
$ws.Cells.Item(67, 13).Value = -4214.778
$ws.Cells.Item(67, 14).Value = -8544.6

$ws.Cells.Item(111, 8).Value = 943.625
$ws.Cells.Item(111, 9).Value = 712.5
$ws.Cells.Item(111, 10).Value = 1174.75
$ws.Cells.Item(111, 11).Value = 2137.5
$ws.Cells.Item(111, 12).Value = 3524.25
$ws.Cells.Item(111, 13).Value = 929.5
$ws.Cells.Item(111, 14).Value = -9658.25

$ws.Cells.Item(132, 8).Value = 1493.5
$ws.Cells.Item(132, 9).Value = 1416.75
$ws.Cells.Item(132, 10).Value = 1723.75
$ws.Cells.Item(132, 11).Value = 4250.25
$ws.Cells.Item(132, 12).Value = 5171.25
$ws.Cells.Item(132, 13).Value = -1720.25
$ws.Cells.Item(132, 14).Value = -10231.25

$ws.Cells.Item(133, 8).Value = 77377.89999999999
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 77377.89999999999
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 77377.89999999999
$ws.Cells.Item(133, 14).Value = -87497.89999999999

$ws.Cells.Item(134, 8).Value = 99995
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 10).Value = 99995
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 12).Value = 99995
$ws.Cells.Item(134, 14).Value = -110135

$ws.Cells.Item(136, 8).Value = 77988.14
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 10).Value = 77988.14
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 12).Value = 77988.14
$ws.Cells.Item(136, 14).Value = -88188.14

$ws.Cells.Item(138, 8).Value = 1494.091
$ws.Cells.Item(138, 9).Value = 1494.091
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 11).Value = 4482.272999999999
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 13).Value = 657.7270000000008
$ws.Cells.Item(138, 14).ClearContents()

$ws.Cells.Item(139, 8).Value = 70767.5
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 70767.5
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 70767.5
$ws.Cells.Item(139, 14).Value = -81047.5

$ws.Cells.Item(140, 8).Value = 50544.11
$ws.Cells.Item(140, 9).Value = 50604.5
$ws.Cells.Item(140, 10).Value = 50526.855
$ws.Cells.Item(140, 11).Value = 50604.5
$ws.Cells.Item(140, 12).Value = 50526.855
$ws.Cells.Item(140, 13).Value = -45424.5
$ws.Cells.Item(140, 14).Value = -60886.855

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 8057.3584
$ws.Cells.Item(32, 9).Value = 3443.2058
$ws.Cells.Item(32, 10).Value = 16314.263
$ws.Cells.Item(32, 11).Value = 3443.2058
$ws.Cells.Item(32, 12).Value = 16314.263
$ws.Cells.Item(32, 13).Value = -3156.2058

$ws.Cells.Item(97, 8).Value = 1619.8
$ws.Cells.Item(97, 9).Value = 1524.75
$ws.Cells.Item(97, 10).Value = 2000
$ws.Cells.Item(97, 11).Value = 1524.75
$ws.Cells.Item(97, 12).Value = 2000
$ws.Cells.Item(97, 13).Value = -1028.75

$ws.Cells.Item(110, 8).Value = 1091.4445
$ws.Cells.Item(110, 9).Value = 671.2143
$ws.Cells.Item(110, 10).Value = 2562.25
$ws.Cells.Item(110, 11).Value = 671.2143
$ws.Cells.Item(110, 12).Value = 2562.25
$ws.Cells.Item(110, 13).Value = 1373.7857
$ws.Cells.Item(110, 14).Value = -6652.25

$ws.Cells.Item(123, 8).Value = 52098.332
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 52098.332
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 12).Value = 52098.332
$ws.Cells.Item(123, 14).Value = -61898.332

$ws.Cells.Item(128, 8).Value = 67700
$ws.Cells.Item(128, 9).Value = 0
$ws.Cells.Item(128, 10).Value = 67700
$ws.Cells.Item(128, 11).Value = 0
$ws.Cells.Item(128, 12).Value = 67700
$ws.Cells.Item(128, 14).Value = -77660

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 1304.1
$ws.Cells.Item(64, 9).Value = 950.25
$ws.Cells.Item(64, 10).Value = 1540
$ws.Cells.Item(64, 11).Value = 950.25
$ws.Cells.Item(64, 12).Value = 1540
$ws.Cells.Item(64, 13).Value = -725.25
$ws.Cells.Item(64, 14).Value = -1990

$ws.Cells.Item(67, 8).Value = 1304.1
$ws.Cells.Item(67, 9).Value = 950.25
$ws.Cells.Item(67, 10).Value = 1540
$ws.Cells.Item(67, 11).Value = 950.25
$ws.Cells.Item(67, 12).Value = 1540
$ws.Cells.Item(67, 13).Value = -170.25
$ws.Cells.Item(67, 14).Value = -3100

$ws.Cells.Item(94, 8).Value = 3445.2727
$ws.Cells.Item(94, 9).Value = 2237.875
$ws.Cells.Item(94, 10).Value = 6665
$ws.Cells.Item(94, 11).Value = 2237.875
$ws.Cells.Item(94, 12).Value = 6665
$ws.Cells.Item(94, 13).Value = -1786.875
$ws.Cells.Item(94, 14).Value = -7567

$ws.Cells.Item(105, 8).Value = 38521.035
$ws.Cells.Item(105, 9).Value = 44579.043
$ws.Cells.Item(105, 10).Value = 3687.5
$ws.Cells.Item(105, 11).Value = 44579.043
$ws.Cells.Item(105, 12).Value = 3687.5
$ws.Cells.Item(105, 13).Value = -42832.043
$ws.Cells.Item(105, 14).Value = -7181.5

$ws.Cells.Item(134, 8).Value = 4293.724
$ws.Cells.Item(134, 9).Value = 2421.3635
$ws.Cells.Item(134, 10).Value = 10178.286
$ws.Cells.Item(134, 11).Value = 7264.0905
$ws.Cells.Item(134, 12).Value = 30534.858
$ws.Cells.Item(134, 13).Value = -4729.0905

$ws.Cells.Item(138, 8).Value = 77893.5
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 77893.5
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 77893.5
$ws.Cells.Item(138, 14).Value = -88173.5

$ws.Cells.Item(140, 8).Value = 93496
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 93496
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 93496
$ws.Cells.Item(140, 14).Value = -103856

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3550
$ws.Cells.Item(31, 9).Value = 2411.625
$ws.Cells.Item(31, 10).Value = 5371.4
$ws.Cells.Item(31, 11).Value = 2411.625
$ws.Cells.Item(31, 12).Value = 5371.4
$ws.Cells.Item(31, 13).Value = -2116.625

$ws.Cells.Item(34, 8).Value = 3550
$ws.Cells.Item(34, 9).Value = 2411.625
$ws.Cells.Item(34, 10).Value = 5371.4
$ws.Cells.Item(34, 11).Value = 2411.625
$ws.Cells.Item(34, 12).Value = 5371.4
$ws.Cells.Item(34, 13).Value = -2209.625

$ws.Cells.Item(69, 8).Value = 44666
$ws.Cells.Item(69, 9).Value = 44666
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 44666
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 13).Value = -43917

$ws.Cells.Item(72, 8).Value = 44666
$ws.Cells.Item(72, 9).Value = 44666
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 11).Value = 133998
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 13).Value = -130254

$ws.Cells.Item(86, 8).Value = 3252689.5
$ws.Cells.Item(86, 9).Value = 5956182
$ws.Cells.Item(86, 10).Value = 8498.799999999999
$ws.Cells.Item(86, 11).Value = 5956182
$ws.Cells.Item(86, 12).Value = 8498.799999999999
$ws.Cells.Item(86, 13).Value = -5955059
$ws.Cells.Item(86, 14).Value = -10744.8

$ws.Cells.Item(89, 8).Value = 3252689.5
$ws.Cells.Item(89, 9).Value = 5956182
$ws.Cells.Item(89, 10).Value = 8498.799999999999
$ws.Cells.Item(89, 11).Value = 29780910
$ws.Cells.Item(89, 12).Value = 42494
$ws.Cells.Item(89, 13).Value = -29775294
$ws.Cells.Item(89, 14).Value = -53726

$ws.Cells.Item(94, 8).Value = 949.25
$ws.Cells.Item(94, 9).Value = 950
$ws.Cells.Item(94, 10).Value = 948.5
$ws.Cells.Item(94, 11).Value = 950
$ws.Cells.Item(94, 12).Value = 948.5
$ws.Cells.Item(94, 13).Value = -499
$ws.Cells.Item(94, 14).Value = -1850.5

$ws.Cells.Item(138, 8).Value = 54918
$ws.Cells.Item(138, 9).Value = 50000
$ws.Cells.Item(138, 10).Value = 55464.445
$ws.Cells.Item(138, 11).Value = 50000
$ws.Cells.Item(138, 12).Value = 55464.445
$ws.Cells.Item(138, 13).Value = -44860
$ws.Cells.Item(138, 14).Value = -65744.44500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(108, 8).Value = 693
$ws.Cells.Item(108, 9).Value = 693
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 11).Value = 2079
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 13).Value = 801

$ws.Cells.Item(109, 8).Value = 3040.9092
$ws.Cells.Item(109, 9).Value = 1362.5
$ws.Cells.Item(109, 10).Value = 4000
$ws.Cells.Item(109, 11).Value = 4087.5
$ws.Cells.Item(109, 12).Value = 12000
$ws.Cells.Item(109, 13).Value = -3047.5

$ws.Cells.Item(110, 8).Value = 7350
$ws.Cells.Item(110, 9).Value = 7400
$ws.Cells.Item(110, 10).Value = 7333.3335
$ws.Cells.Item(110, 11).Value = 22200
$ws.Cells.Item(110, 12).Value = 22000.0005
$ws.Cells.Item(110, 13).Value = -18110

$ws.Cells.Item(112, 8).Value = 5521.393
$ws.Cells.Item(112, 9).Value = 4499.5
$ws.Cells.Item(112, 10).Value = 5600
$ws.Cells.Item(112, 11).Value = 13498.5
$ws.Cells.Item(112, 12).Value = 16800
$ws.Cells.Item(112, 13).Value = -12390.5

$ws.Cells.Item(116, 8).Value = 2164.8
$ws.Cells.Item(116, 9).Value = 2206
$ws.Cells.Item(116, 10).Value = 2000
$ws.Cells.Item(116, 11).Value = 6618
$ws.Cells.Item(116, 12).Value = 6000
$ws.Cells.Item(116, 13).Value = -3176

$ws.Cells.Item(129, 8).Value = 41667140
$ws.Cells.Item(129, 9).Value = 545.3333
$ws.Cells.Item(129, 10).Value = 166666930
$ws.Cells.Item(129, 11).Value = 1635.9999
$ws.Cells.Item(129, 12).Value = 500000790
$ws.Cells.Item(129, 13).Value = 3364.0001
$ws.Cells.Item(129, 14).Value = -500010790

$ws.Cells.Item(132, 8).Value = 6310.5
$ws.Cells.Item(132, 9).Value = 1536
$ws.Cells.Item(132, 10).Value = 7265.4
$ws.Cells.Item(132, 11).Value = 13824
$ws.Cells.Item(132, 12).Value = 65388.6
$ws.Cells.Item(132, 13).Value = -11294
$ws.Cells.Item(132, 14).Value = -70448.60000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(140, 8).Value = 94518
$ws.Cells.Item(140, 9).Value = 90500
$ws.Cells.Item(140, 10).Value = 94964.44500000001
$ws.Cells.Item(140, 11).Value = 90500
$ws.Cells.Item(140, 12).Value = 94964.44500000001
$ws.Cells.Item(140, 13).Value = -85320
$ws.Cells.Item(140, 14).Value = -105324.445

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1609.0714
$ws.Cells.Item(46, 9).Value = 1594.3846
$ws.Cells.Item(46, 10).Value = 1800
$ws.Cells.Item(46, 11).Value = 1594.3846
$ws.Cells.Item(46, 12).Value = 1800
$ws.Cells.Item(46, 13).Value = -1406.3846

$ws.Cells.Item(55, 8).Value = 4878689.5
$ws.Cells.Item(55, 9).Value = 493.4
$ws.Cells.Item(55, 10).Value = 7693033.5
$ws.Cells.Item(55, 11).Value = 493.4
$ws.Cells.Item(55, 12).Value = 7693033.5
$ws.Cells.Item(55, 13).Value = -320.4
$ws.Cells.Item(55, 14).Value = -7693379.5

$ws.Cells.Item(93, 8).Value = 1240.5555
$ws.Cells.Item(93, 9).Value = 1022.5714
$ws.Cells.Item(93, 10).Value = 2003.5
$ws.Cells.Item(93, 11).Value = 1022.5714
$ws.Cells.Item(93, 12).Value = 2003.5
$ws.Cells.Item(93, 13).Value = 225.4286
$ws.Cells.Item(93, 14).Value = -4499.5

$ws.Cells.Item(114, 8).Value = 0
$ws.Cells.Item(114, 9).Value = 0
$ws.Cells.Item(114, 10).Value = 0
$ws.Cells.Item(114, 11).Value = 0
$ws.Cells.Item(114, 12).Value = 0
$ws.Cells.Item(114, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 1674
$ws.Cells.Item(132, 9).Value = 1732.6666
$ws.Cells.Item(132, 10).Value = 1498
$ws.Cells.Item(132, 11).Value = 5197.9998
$ws.Cells.Item(132, 12).Value = 4494
$ws.Cells.Item(132, 13).Value = -2667.9998
$ws.Cells.Item(132, 14).Value = -9554

$ws.Cells.Item(136, 8).Value = 4026.7273
$ws.Cells.Item(136, 9).Value = 4389.8667
$ws.Cells.Item(136, 10).Value = 3248.5715
$ws.Cells.Item(136, 11).Value = 13169.6001
$ws.Cells.Item(136, 12).Value = 9745.7145
$ws.Cells.Item(136, 13).Value = -10619.6001
